$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text (string) valued cells - force text format so Excel does not
# re-interpret the numeric-looking strings as numbers.
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "127.42000000"

$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "85596.04072000"

$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "10903507.61611180"

$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "43752.55114000"

$ws.Range("K20").NumberFormat = "@"
$ws.Range("K20").Value = "5574866.68277130"

# Numeric valued cells
$ws.Range("I20").Value = 38744

$ws.Range("M20").Value = 127.4199999999999
$ws.Range("N20").Value = 126.855
$ws.Range("O20").Value = 127.4314285714286
$ws.Range("P20").Value = 129.796
$ws.Range("R20").Value = 127.42
$ws.Range("S20").Value = 126.9831169962921
$ws.Range("T20").Value = 128.5404561014247
$ws.Range("U20").Value = 130.4611056282133
$ws.Range("V20").Value = -1.920649526788651
$ws.Range("W20").Value = -1.807482571254994
$ws.Range("X20").Value = -0.1131669555336565
